# Update crypto price/volume figures per the scraped refresh (2024-04-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.352.15'
$ws.Range("E2").Value = '  -2.58%  '

$ws.Range("D3").Value = '3.534.10'
$ws.Range("E3").Value = '  -4.36%  '

$ws.Range("E4").Value = '  -0.07%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'580.86"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -0.14%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'171.97"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -3.52%  '

$ws.Range("D7").Value = '3.530.80'
$ws.Range("E7").Value = '  -4.11%  '

$ws.Range("E8").Value = '  -1.17%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("E10").Value = '  -5.02%  '

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'6.66"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -3.81%  '

$ws.Range("E12").Value = '  -4.02%  '

$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'47.46"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -3.52%  '

$ws.Range("E14").Value = '  -4.63%  '

$ws.Range("D15").Value = '4.094.05'
$ws.Range("E15").Value = '  -4.67%  '

$ws.Range("E16").Value = '  -5.03%  '

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'629.12"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -7.27%  '

$ws.Range("D18").Value = '3.533.73'
$ws.Range("E18").Value = '  -4.42%  '

$ws.Range("D19").Value = '69.298.45'
$ws.Range("E19").Value = '  -2.86%  '

$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("E21").Value = '  -2.52%  '

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'11.23"
$ws.Range("D22").Style = $style

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.891"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -5.64%  '

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'16.04"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -7.89%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'97.92"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -4.08%  '

$ws.Range("E26").Value = '  -4.24%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("E28").Value = '  -6.91%  '

$ws.Range("E29").Value = '  -9.18%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'32.93"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -6.25%  '

$ws.Range("E31").Value = '  -7.64%  '

$ws.Range("E32").Value = '  -5.98%  '

$ws.Range("E33").Value = '  -6.62%  '

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'7.03"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -6.56%  '

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'631.93"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  +8.77%  '

$ws.Range("E36").Value = '  -3.66%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'3.51"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -13.07%  '

$ws.Range("E38").Value = '  -4.57%  '

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'57.39"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -2.33%  '

$ws.Range("E40").Value = '  +0.00%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.0457"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -1.54%  '

$ws.Range("E42").Value = '  -5.40%  '

$ws.Range("D43").Value = '3.398.70'
$ws.Range("E43").Value = '  -6.30%  '

$ws.Range("E44").Value = '  -6.15%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'33.06"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -7.11%  '

$ws.Range("D46").Value = '0.0₃0700'
$ws.Range("E46").Value = '  -8.96%  '

$ws.Range("E47").Value = '  -7.04%  '

$ws.Range("E48").Value = '  -4.40%  '

$ws.Range("E49").Value = '  -2.46%  '

$ws.Range("E50").Value = '  +14.74%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'131.59"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -2.26%  '
